$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2400
$ws.Range("J28").Value = 1500
$ws.Range("L28").Value = 1500
$ws.Range("N28").Value = -2470
$ws.Range("H33").Value = 367.72726
$ws.Range("I33").Value = 222.57143
$ws.Range("K33").Value = 222.57143
$ws.Range("M33").Value = 6.428570000000008
$ws.Range("H62").Value = 9342.571
$ws.Range("I62").Value = 11720
$ws.Range("K62").Value = 11720
$ws.Range("M62").Value = -11096
$ws.Range("H65").Value = 9342.571
$ws.Range("I65").Value = 11720
$ws.Range("K65").Value = 58600
$ws.Range("M65").Value = -55480
$ws.Range("H92").Value = 612.94116
$ws.Range("I92").Value = 137.14285
$ws.Range("J92").Value = 2833.3333
$ws.Range("K92").Value = 137.14285
$ws.Range("L92").Value = 2833.3333
$ws.Range("M92").Value = 1110.85715
$ws.Range("N92").Value = -5329.3333
$ws.Range("H100").Value = 2299.8
$ws.Range("J100").Value = 2999.5
$ws.Range("L100").Value = 2999.5
$ws.Range("N100").Value = -4081.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4900
$ws.Range("I45").Value = 4875
$ws.Range("K45").Value = 4875
$ws.Range("M45").Value = -4498
$ws.Range("H74").Value = 16220.071
$ws.Range("I74").Value = 2097.0908
$ws.Range("K74").Value = 2097.0908
$ws.Range("M74").Value = -1223.0908
$ws.Range("H77").Value = 16220.071
$ws.Range("I77").Value = 2097.0908
$ws.Range("K77").Value = 10485.454
$ws.Range("M77").Value = -6117.454
$ws.Range("H122").Value = 2174.2104
$ws.Range("I122").Value = 1986
$ws.Range("J122").Value = 2582
$ws.Range("K122").Value = 5958
$ws.Range("L122").Value = 7746
$ws.Range("M122").Value = -3508
$ws.Range("N122").Value = -12646
$ws.Range("H132").Value = 1390134.6
$ws.Range("I132").Value = 1471847.8
$ws.Range("K132").Value = 4415543.4
$ws.Range("M132").Value = -4413013.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("N39").Value = -10778
$ws.Range("H105").Value = 2824.5112
$ws.Range("I105").Value = 2190.4194
$ws.Range("J105").Value = 4228.5713
$ws.Range("K105").Value = 2190.4194
$ws.Range("L105").Value = 4228.5713
$ws.Range("M105").Value = -443.4194000000002
$ws.Range("N105").Value = -7722.5713
$ws.Range("H107").Value = 1086.5
$ws.Range("I107").Value = 1070.2858
$ws.Range("K107").Value = 1070.2858
$ws.Range("M107").Value = 849.7141999999999
$ws.Range("H134").Value = 11429.143
$ws.Range("I134").Value = 4000.8
$ws.Range("K134").Value = 12002.4
$ws.Range("M134").Value = -9467.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 870.2143
$ws.Range("I16").Value = 833.9091
$ws.Range("J16").Value = 1003.3333
$ws.Range("K16").Value = 833.9091
$ws.Range("L16").Value = 1003.3333
$ws.Range("M16").Value = -546.9091
$ws.Range("N16").Value = -1577.3333
$ws.Range("H99").Value = 3979622.2
$ws.Range("I99").Value = 27543.295
$ws.Range("K99").Value = 27543.295
$ws.Range("M99").Value = -26045.295
$ws.Range("H105").Value = 14767
$ws.Range("J105").Value = 4899.5
$ws.Range("L105").Value = 4899.5
$ws.Range("N105").Value = -8393.5
$ws.Range("H113").Value = 870.2143
$ws.Range("I113").Value = 833.9091
$ws.Range("J113").Value = 1003.3333
$ws.Range("K113").Value = 833.9091
$ws.Range("L113").Value = 1003.3333
$ws.Range("M113").Value = 1336.0909
$ws.Range("N113").Value = -5343.3333
$ws.Range("H126").Value = 3979622.2
$ws.Range("I126").Value = 27543.295
$ws.Range("K126").Value = 82629.88499999999
$ws.Range("M126").Value = -80159.88499999999
$ws.Range("H134").Value = 4120.2666
$ws.Range("I134").Value = 3058.75
$ws.Range("K134").Value = 9176.25
$ws.Range("M134").Value = -6641.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 240.41667
$ws.Range("I60").Value = 208.7
$ws.Range("K60").Value = 626.0999999999999
$ws.Range("M60").Value = -375.0999999999999
$ws.Range("H107").Value = 2107.8333
$ws.Range("I107").Value = 550
$ws.Range("K107").Value = 1650
$ws.Range("M107").Value = 270
$ws.Range("H129").Value = 1790.1765
$ws.Range("I129").Value = 943.125
$ws.Range("J129").Value = 2543.111
$ws.Range("K129").Value = 2829.375
$ws.Range("L129").Value = 7629.333
$ws.Range("M129").Value = 2170.625
$ws.Range("N129").Value = -17629.333
$ws.Range("H137").Value = 7518.909
$ws.Range("J137").Value = 12379.8
$ws.Range("L137").Value = 37139.39999999999
$ws.Range("N137").Value = -47339.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1539.289
$ws.Range("I102").Value = 1702.7435
$ws.Range("J102").Value = 476.83334
$ws.Range("K102").Value = 1702.7435
$ws.Range("L102").Value = 476.83334
$ws.Range("M102").Value = -80.74350000000004
$ws.Range("N102").Value = -3720.83334
$ws.Range("H113").Value = 3206.2
$ws.Range("I113").Value = 2216.2
$ws.Range("J113").Value = 6176.2
$ws.Range("K113").Value = 2216.2
$ws.Range("L113").Value = 6176.2
$ws.Range("M113").Value = -46.19999999999982
$ws.Range("N113").Value = -10516.2
$ws.Range("H122").Value = 3221.16
$ws.Range("I122").Value = 2767.0588
$ws.Range("J122").Value = 4186.125
$ws.Range("K122").Value = 8301.1764
$ws.Range("L122").Value = 12558.375
$ws.Range("M122").Value = -5851.1764
$ws.Range("N122").Value = -17458.375
$ws.Range("H123").Value = 49517.332
$ws.Range("J123").Value = 49517.332
$ws.Range("L123").Value = 49517.332
$ws.Range("N123").Value = -54417.332
$ws.Range("H126").Value = 6194.0586
$ws.Range("I126").Value = 6885.643
$ws.Range("K126").Value = 20656.929
$ws.Range("M126").Value = -18186.929
$ws.Range("H132").Value = 20290.545
$ws.Range("I132").Value = 22723.555
$ws.Range("J132").Value = 9342
$ws.Range("K132").Value = 68170.66500000001
$ws.Range("L132").Value = 28026
$ws.Range("M132").Value = -65640.66500000001
$ws.Range("N132").Value = -33086

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3594.9375
$ws.Range("J7").Value = 3214.8333
$ws.Range("L7").Value = 3214.8333
$ws.Range("N7").Value = -3438.8333
$ws.Range("H40").Value = 4139.222
$ws.Range("I40").Value = 3813.25
$ws.Range("K40").Value = 3813.25
$ws.Range("M40").Value = -3677.25
$ws.Range("H46").Value = 3023.6316
$ws.Range("I46").Value = 800.6
$ws.Range("K46").Value = 800.6
$ws.Range("M46").Value = -612.6
$ws.Range("H61").Value = 23299.75
$ws.Range("I61").Value = 21057.143
$ws.Range("J61").Value = 38998
$ws.Range("K61").Value = 21057.143
$ws.Range("L61").Value = 38998
$ws.Range("M61").Value = -20855.143
$ws.Range("N61").Value = -39402
$ws.Range("H82").Value = 1699.5834
$ws.Range("I82").Value = 3633
$ws.Range("K82").Value = 3633
$ws.Range("M82").Value = -3272
$ws.Range("H85").Value = 1699.5834
$ws.Range("I85").Value = 3633
$ws.Range("K85").Value = 3633
$ws.Range("M85").Value = -2385
$ws.Range("H100").Value = 3742.7646
$ws.Range("I100").Value = 3495.1333
$ws.Range("K100").Value = 3495.1333
$ws.Range("M100").Value = -2954.1333
$ws.Range("H113").Value = 23299.75
$ws.Range("I113").Value = 21057.143
$ws.Range("J113").Value = 38998
$ws.Range("K113").Value = 21057.143
$ws.Range("L113").Value = 38998
$ws.Range("M113").Value = -18887.143
$ws.Range("N113").Value = -43338
$ws.Range("H126").Value = 3594.9375
$ws.Range("J126").Value = 3214.8333
$ws.Range("L126").Value = 9644.499899999999
$ws.Range("N126").Value = -14584.4999
$ws.Range("H136").Value = 10066.866
$ws.Range("I136").Value = 3989.6667
$ws.Range("J136").Value = 11586.167
$ws.Range("K136").Value = 11969.0001
$ws.Range("L136").Value = 34758.501
$ws.Range("M136").Value = -9419.000100000001
$ws.Range("N136").Value = -39858.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2044
$ws.Range("I107").Value = 1190.1578
$ws.Range("K107").Value = 3570.4734
$ws.Range("M107").Value = -1650.4734
$ws.Range("H126").Value = 2283.5881
$ws.Range("I126").Value = 2356
$ws.Range("K126").Value = 7068
$ws.Range("M126").Value = -4598
$ws.Range("H132").Value = 2120.95
$ws.Range("I132").Value = 1949.1765
$ws.Range("J132").Value = 3094.3333
$ws.Range("K132").Value = 5847.529500000001
$ws.Range("L132").Value = 9282.999899999999
$ws.Range("M132").Value = -3317.529500000001
$ws.Range("N132").Value = -14342.9999
$ws.Range("H137").Value = 53666.668
$ws.Range("I137").Value = 46000
$ws.Range("J137").Value = 57500
$ws.Range("K137").Value = 46000
$ws.Range("L137").Value = 57500
$ws.Range("M137").Value = -40900
$ws.Range("N137").Value = -67700
